# Journal de travail - mise a jour du code et de la documentation
# Applies the content changes for rows 30 and 31 (Feuil1):
#   - fills in date/heure/module/type/tache/emplacement/descriptif for row 30
#   - fills in date/heure/module/type for row 31, with an anomalous text value
#     in the "Heure fin" column (D31) that turns the shared duration formula
#     into a #VALUE! error, and a hyperlinked source in K31
#   - updates the active selection / scroll position to match the new editing
#     location

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------
# Row 30 : "Finalisation de la Bataille Navale 0.1"
# ---------------------------------------------------------------------
$ws.Range("B30").Value = 44267
$ws.Range("C30").Value = 42217.510416666664
$ws.Range("D30").Value = 42217.5625
$ws.Range("F30").Value = "MA-20"
$ws.Range("H30").Value = "Finalisation de la Bataille Navale 0.1"
$ws.Range("J30").Value = "Avec quelques modifications sur mon code, j'ai réussi à faire gagner.`nJ'ai enlevé un if qui posait probleme et inversé le ""do… while"" qui ignorait les variables à la fin."

# Row 30 grows to fit the new multi-line description.
$ws.Rows.Item(30).RowHeight = 73

# ---------------------------------------------------------------------
# Row 31 : follow-up entry, with a malformed "Heure fin" (D31) entered as
# free text instead of a time, which breaks the shared duration formula.
# ---------------------------------------------------------------------
$ws.Range("B31").Value = 44268
$ws.Range("C31").Value = 42217.5625
$ws.Range("F31").Value = "MA-20"

$ws.Hyperlinks.Add($ws.Range("K31"), "https://cboard.cprogramming.com/cplusplus-programming/117049-setconsoletitle.html") | Out-Null

$ws.Range("D31").Value = "01.08.2015  :00:00"

# ---------------------------------------------------------------------
# Selection / scroll position left where the author was working
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H31").Select() | Out-Null
